$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -7.3854454316297193
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -7.1311412664558205

$ws.Range("B3").Value = -8.6193128364008444
$ws.Range("C3").Value = -3.4892700104559182
$ws.Range("D3").Value = -10.505396392868107
$ws.Range("E3").Value = 8.3121526468800937

$ws.Range("B1:E3").Select() | Out-Null
